$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.440.36"
$ws.Range("E2").Value = "  +2.62%  "

$ws.Range("D3").Value = "3.143.22"
$ws.Range("E3").Value = "  +2.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.69"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.28"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +6.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "3.142.15"
$ws.Range("E8").Value = "  +2.28%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.52"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.93%  "

$ws.Range("E11").Value = "  +1.88%  "

$ws.Range("E12").Value = "  +0.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000243"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.31%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.99"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.68%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.667.38"
$ws.Range("E15").Value = "  +2.31%  "

$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.122"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.94%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "68.340.69"
$ws.Range("E17").Value = "  +2.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.15"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.35%  "

$ws.Range("D19").Value = "3.141.38"
$ws.Range("E19").Value = "  +2.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.54"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.86%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "488.83"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.16%  "

$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.699"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.59%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.82"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.00"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.57%  "

$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +6.38%  "

$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.02"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.77%  "

$ws.Range("E27").Value = "  +4.24%  "

$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.13"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.37"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.72%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.65"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.29"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.61%  "

$ws.Range("E33").Value = "  +1.28%  "

$ws.Range("D34").Value = "0.0₃0952"
$ws.Range("E34").Value = "  +4.09%  "

$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "48.86"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.69"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.81%  "

$ws.Range("E38").Value = "  +1.47%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.326"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +8.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.06"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "49.22"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.125"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.17%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.42"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.44%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.75"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +9.05%  "

$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "403.08"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +9.39%  "

$ws.Range("D46").Value = "2.804.02"
$ws.Range("E46").Value = "  +1.51%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.41"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +11.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0350"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.49%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "135.60"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.36%  "

$ws.Range("E50").Value = "  +0.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.37"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +10.13%  "
